# "Added more struggle data and repaired falling data"
# The "falling" sample rows (2-21, columns C:H = ax,ay,az,gx,gy,gz) are
# shifted up by one reading: each row now holds the sensor reading that
# used to belong to the row below it, and the final row (21) is populated
# with one freshly captured reading. Timestamps (col A) and the "falling"
# label (col B) are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = -0.2705469131469726
$ws.Cells.Item(2, 4).Value = 0.4402385950088501
$ws.Cells.Item(2, 5).Value = -1.731513738632202
$ws.Cells.Item(2, 6).Value = -0.05105815259351689
$ws.Cells.Item(2, 7).Value = -0.02509637922048548
$ws.Cells.Item(2, 8).Value = -0.0663297846913338

$ws.Cells.Item(3, 3).Value = -0.6277971267700195
$ws.Cells.Item(3, 4).Value = 0.0076048374176025
$ws.Cells.Item(3, 5).Value = -1.507715225219727
$ws.Cells.Item(3, 6).Value = -0.1197514058578581
$ws.Cells.Item(3, 7).Value = -0.3090105539276482
$ws.Cells.Item(3, 8).Value = 0.06041020864532091

$ws.Cells.Item(4, 3).Value = -2.914698600769043
$ws.Cells.Item(4, 4).Value = -1.449564576148987
$ws.Cells.Item(4, 5).Value = -3.32840347290039
$ws.Cells.Item(4, 6).Value = -0.4005748778581615
$ws.Cells.Item(4, 7).Value = -0.7805985297475537
$ws.Cells.Item(4, 8).Value = 0.1577123148100716

$ws.Cells.Item(5, 3).Value = 1.028462886810303
$ws.Cells.Item(5, 4).Value = -0.5502710342407227
$ws.Cells.Item(5, 5).Value = -4.842555046081543
$ws.Cells.Item(5, 6).Value = -0.6522004490806942
$ws.Cells.Item(5, 7).Value = -1.329693669364566
$ws.Cells.Item(5, 8).Value = 0.1308778794038864

$ws.Cells.Item(6, 3).Value = -1.341280460357666
$ws.Cells.Item(6, 4).Value = -2.225003957748413
$ws.Cells.Item(6, 5).Value = -6.344600677490234
$ws.Cells.Item(6, 6).Value = -0.6102398293358944
$ws.Cells.Item(6, 7).Value = -1.220246967815219
$ws.Cells.Item(6, 8).Value = 0.4809618578070668

$ws.Cells.Item(7, 3).Value = 8.574896812438965
$ws.Cells.Item(7, 4).Value = 0.6133027076721191
$ws.Cells.Item(7, 5).Value = -6.888121604919434
$ws.Cells.Item(7, 6).Value = -0.04692753723689685
$ws.Cells.Item(7, 7).Value = -0.2711587122508481
$ws.Cells.Item(7, 8).Value = 1.261436768940515

$ws.Cells.Item(8, 3).Value = -6.096681118011475
$ws.Cells.Item(8, 4).Value = 0.8472604751586914
$ws.Cells.Item(8, 5).Value = 14.72706890106201
$ws.Cells.Item(8, 6).Value = 0.2368920927955998
$ws.Cells.Item(8, 7).Value = 2.034886604263666
$ws.Cells.Item(8, 8).Value = 0.6732607796078658

$ws.Cells.Item(9, 3).Value = 4.274323463439941
$ws.Cells.Item(9, 4).Value = -4.468049049377441
$ws.Cells.Item(9, 5).Value = -6.856836795806885
$ws.Cells.Item(9, 6).Value = -0.7203119397163371
$ws.Cells.Item(9, 7).Value = 3.954537868499759
$ws.Cells.Item(9, 8).Value = -2.220546166102086

$ws.Cells.Item(10, 3).Value = -4.518700122833252
$ws.Cells.Item(10, 4).Value = -1.648021101951599
$ws.Cells.Item(10, 5).Value = -0.9248533248901368
$ws.Cells.Item(10, 6).Value = -0.9733701603753284
$ws.Cells.Item(10, 7).Value = 2.440581185477122
$ws.Cells.Item(10, 8).Value = -3.891320841653013

$ws.Cells.Item(11, 3).Value = 9.755411148071287
$ws.Cells.Item(11, 4).Value = 3.367114305496216
$ws.Cells.Item(11, 5).Value = 2.822277307510376
$ws.Cells.Item(11, 6).Value = 0.4988514525549754
$ws.Cells.Item(11, 7).Value = 0.8867653551555867
$ws.Cells.Item(11, 8).Value = -1.71078631139937

$ws.Cells.Item(12, 3).Value = 1.561064720153809
$ws.Cells.Item(12, 4).Value = 0.1129603385925293
$ws.Cells.Item(12, 5).Value = -0.9029455184936525
$ws.Cells.Item(12, 6).Value = 0.1264272814705268
$ws.Cells.Item(12, 7).Value = -2.447191684019
$ws.Cells.Item(12, 8).Value = -0.4622068021978656

$ws.Cells.Item(13, 3).Value = 5.92741584777832
$ws.Cells.Item(13, 4).Value = -0.8555939197540283
$ws.Cells.Item(13, 5).Value = 4.797466278076172
$ws.Cells.Item(13, 6).Value = -0.7295694393771054
$ws.Cells.Item(13, 7).Value = -0.004014266388761123
$ws.Cells.Item(13, 8).Value = 0.09374600010258789

$ws.Cells.Item(14, 3).Value = 1.122594833374023
$ws.Cells.Item(14, 4).Value = 1.295500755310059
$ws.Cells.Item(14, 5).Value = -1.442571401596069
$ws.Cells.Item(14, 6).Value = -0.1179624412740978
$ws.Cells.Item(14, 7).Value = 1.259269575277969
$ws.Cells.Item(14, 8).Value = 0.1631955632141657

$ws.Cells.Item(15, 3).Value = 0.5986118316650391
$ws.Cells.Item(15, 4).Value = 0.4096674025058746
$ws.Cells.Item(15, 5).Value = -0.6679027080535889
$ws.Cells.Item(15, 6).Value = -0.1617193005624273
$ws.Cells.Item(15, 7).Value = -0.4078179995218952
$ws.Cells.Item(15, 8).Value = 0.2281727109636559

$ws.Cells.Item(16, 3).Value = 0.0388402938842773
$ws.Cells.Item(16, 4).Value = 0.3524296283721924
$ws.Cells.Item(16, 5).Value = -1.101761341094971
$ws.Cells.Item(16, 6).Value = -0.124921940267086
$ws.Cells.Item(16, 7).Value = 1.249372124671936
$ws.Cells.Item(16, 8).Value = 1.016479730606079

$ws.Cells.Item(17, 3).Value = -0.1728830337524414
$ws.Cells.Item(17, 4).Value = 0.6193998456001282
$ws.Cells.Item(17, 5).Value = -0.6873818635940552
$ws.Cells.Item(17, 6).Value = -0.4928955077415405
$ws.Cells.Item(17, 7).Value = -0.4476696934018742
$ws.Cells.Item(17, 8).Value = -0.9833766732896989

$ws.Cells.Item(18, 3).Value = 0.4876585006713867
$ws.Cells.Item(18, 4).Value = 0.6636635065078735
$ws.Cells.Item(18, 5).Value = -0.9166454076766968
$ws.Cells.Item(18, 6).Value = -0.2805471434479678
$ws.Cells.Item(18, 7).Value = 0.4230750912711692
$ws.Cells.Item(18, 8).Value = -0.2188279224293611

$ws.Cells.Item(19, 3).Value = -0.1092472076416015
$ws.Cells.Item(19, 4).Value = 0.732629120349884
$ws.Cells.Item(19, 5).Value = -1.016466021537781
$ws.Cells.Item(19, 6).Value = -0.2866266923291341
$ws.Cells.Item(19, 7).Value = 0.4182899764605935
$ws.Cells.Item(19, 8).Value = -0.003992439912898826

$ws.Cells.Item(20, 3).Value = 0.4153709411621094
$ws.Cells.Item(20, 4).Value = 0.5096800327301025
$ws.Cells.Item(20, 5).Value = -0.7671611309051514
$ws.Cells.Item(20, 6).Value = -0.03713915026968551
$ws.Cells.Item(20, 7).Value = 0.07855436143775912
$ws.Cells.Item(20, 8).Value = 0.07685266648020034

$ws.Cells.Item(21, 3).Value = 0.17730712890625
$ws.Cells.Item(21, 4).Value = 0.6253083348274231
$ws.Cells.Item(21, 5).Value = -0.8837988376617432
$ws.Cells.Item(21, 6).Value = 0.03008511281084441
$ws.Cells.Item(21, 7).Value = 0.0510363349070151
$ws.Cells.Item(21, 8).Value = -0.07897615255344478
